$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.109.52"
$ws.Range("E2").Value = "  -2.13%  "
$ws.Range("D3").Value = "3.477.48"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'587.76"
$ws.Range("E5").Value = "  -3.05%  "
$ws.Range("D6").Value = "'138.16"
$ws.Range("E6").Value = "  -3.12%  "
$ws.Range("D7").Value = "3.474.77"
$ws.Range("E7").Value = "  -1.04%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.488"
$ws.Range("E9").Value = "  -4.06%  "
$ws.Range("E10").Value = "  -5.49%  "
$ws.Range("D11").Value = "'7.19"
$ws.Range("E11").Value = "  -7.01%  "
$ws.Range("D12").Value = "'0.380"
$ws.Range("E12").Value = "  -6.45%  "
$ws.Range("D13").Value = "4.052.31"
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("D14").Value = "'0.0000182"
$ws.Range("E14").Value = "  -5.85%  "
$ws.Range("D15").Value = "'26.57"
$ws.Range("E15").Value = "  -6.84%  "
$ws.Range("D16").Value = "3.453.10"
$ws.Range("E16").Value = "  -1.94%  "
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("D18").Value = "64.898.88"
$ws.Range("E18").Value = "  -2.16%  "
$ws.Range("E19").Value = "  -10.38%  "
$ws.Range("D20").Value = "'5.75"
$ws.Range("E20").Value = "  -6.14%  "
$ws.Range("D21").Value = "'13.79"
$ws.Range("E21").Value = "  -5.79%  "
$ws.Range("D22").Value = "'387.51"
$ws.Range("E22").Value = "  -7.91%  "
$ws.Range("D23").Value = "'0.553"
$ws.Range("E23").Value = "  -5.53%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "'72.39"
$ws.Range("E25").Value = "  -5.85%  "
$ws.Range("D26").Value = "'5.75"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "3.602.12"
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("D28").Value = "'0.0000109"
$ws.Range("E28").Value = "  -3.01%  "
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").Value = "'7.42"
$ws.Range("E30").Value = "  -5.95%  "
$ws.Range("D31").Value = "'8.17"
$ws.Range("E31").Value = "  -8.08%  "
$ws.Range("D32").Value = "'2.21"
$ws.Range("E32").Value = "  -9.78%  "
$ws.Range("D33").Value = "3.479.81"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.143"
$ws.Range("E35").Value = "  -6.95%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "'23.02"
$ws.Range("E36").Value = "  -4.52%  "
$ws.Range("D37").Value = "'171.22"
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("D38").Value = "'1.20"
$ws.Range("E38").Value = "  -9.68%  "
$ws.Range("D39").Value = "'6.87"
$ws.Range("E39").Value = "  -8.88%  "
$ws.Range("D40").Value = "'1.49"
$ws.Range("E40").Value = "  -8.07%  "
$ws.Range("D41").Value = "'4.73"
$ws.Range("E41").Value = "  -8.71%  "
$ws.Range("D42").Value = "'0.0773"
$ws.Range("E42").Value = "  -5.09%  "
$ws.Range("D43").Value = "'0.809"
$ws.Range("E43").Value = "  -4.86%  "
$ws.Range("D44").Value = "'0.998"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "'42.24"
$ws.Range("E45").Value = "  -7.01%  "
$ws.Range("D46").Value = "'4.33"
$ws.Range("E46").Value = "  -12.50%  "
$ws.Range("D47").Value = "'1.62"
$ws.Range("E47").Value = "  -7.85%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'1.13"
$ws.Range("E48").Value = "  +2.70%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'23.32"
$ws.Range("E49").Value = "  +3.35%  "
$ws.Range("D50").Value = "'6.63"
$ws.Range("E50").Value = "  -5.81%  "
$ws.Range("D51").Value = "2.219.76"
$ws.Range("E51").Value = "  -3.66%  "
